$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the second test row to a new workflow file and refresh the run statistics
$ws.Range("A3").Value = "Test_Framework\Tests\_wbTest_Example2.xaml"
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 20

# Move the active selection to B4, matching the post-edit cursor position
$ws.Range("B4").Select()
